$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Designs")

# Insert a new column before column A to hold the design ID
$ws.Columns.Item(1).Insert() | Out-Null

# Header
$ws.Range("A1").Value = "ID"

# ID values for each design row (rows 2-9)
$ws.Range("A2").Value = 1025
$ws.Range("A3").Value = 1009
$ws.Range("A4").Value = 1006
$ws.Range("A5").Value = 1003
$ws.Range("A6").Value = 1
$ws.Range("A7").Value = 11
$ws.Range("A8").Value = 1000
$ws.Range("A9").Value = 1001

# Match the author's final cursor position in the sheet
$ws.Range("Q6").Select() | Out-Null
